$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2750.5
$ws.Range("I94").Value = 2750.5
$ws.Range("K94").Value = 2750.5
$ws.Range("M94").Value = -2299.5
$ws.Range("H98").Value = 58825784
$ws.Range("I98").Value = 62502210
$ws.Range("J98").Value = 2999
$ws.Range("K98").Value = 62502210
$ws.Range("L98").Value = 2999
$ws.Range("M98").Value = -62500712
$ws.Range("N98").Value = -5995
$ws.Range("H100").Value = 2644.353
$ws.Range("I100").Value = 1189.3334
$ws.Range("K100").Value = 1189.3334
$ws.Range("M100").Value = -648.3334
$ws.Range("H112").Value = 2572.2222
$ws.Range("I112").Value = 5000
$ws.Range("J112").Value = 2429.4119
$ws.Range("K112").Value = 15000
$ws.Range("L112").Value = 7288.2357
$ws.Range("M112").Value = -13892
$ws.Range("N112").Value = -9504.235700000001
$ws.Range("H122").Value = 58825784
$ws.Range("I122").Value = 62502210
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 187506630
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -187504180
$ws.Range("N122").Value = -13897
$ws.Range("H132").Value = 1612.2858
$ws.Range("I132").Value = 1560.9259
$ws.Range("K132").Value = 4682.7777
$ws.Range("M132").Value = -2152.7777
$ws.Range("H137").Value = 6971.6924
$ws.Range("I137").Value = 4198.909
$ws.Range("K137").Value = 12596.727
$ws.Range("M137").Value = -10046.727

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11113411
$ws.Range("I32").Value = 12501860
$ws.Range("J32").Value = 5816.2
$ws.Range("K32").Value = 12501860
$ws.Range("L32").Value = 5816.2
$ws.Range("M32").Value = -12501573
$ws.Range("N32").Value = -6390.2
$ws.Range("H38").Value = 2507500
$ws.Range("I38").Value = 2507500
$ws.Range("K38").Value = 2507500
$ws.Range("M38").Value = -2507033
$ws.Range("H61").Value = 21787942
$ws.Range("I61").Value = 45457850
$ws.Range("J61").Value = 90529.336
$ws.Range("K61").Value = 45457850
$ws.Range("L61").Value = 90529.336
$ws.Range("M61").Value = -45457638
$ws.Range("N61").Value = -90953.336
$ws.Range("H102").Value = 8800.799999999999
$ws.Range("I102").Value = 8800.799999999999
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 8800.799999999999
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -7178.799999999999
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value = 1677.9
$ws.Range("I110").Value = 1476.375
$ws.Range("J110").Value = 2484
$ws.Range("K110").Value = 1476.375
$ws.Range("L110").Value = 2484
$ws.Range("M110").Value = 568.625
$ws.Range("N110").Value = -6574
$ws.Range("H132").Value = 4290.7393
$ws.Range("I132").Value = 688.875
$ws.Range("J132").Value = 12523.571
$ws.Range("K132").Value = 2066.625
$ws.Range("L132").Value = 37570.713
$ws.Range("M132").Value = 463.375
$ws.Range("N132").Value = -42630.713
$ws.Range("H134").Value = 64749.75
$ws.Range("J134").Value = 64749.75
$ws.Range("L134").Value = 64749.75
$ws.Range("N134").Value = -74889.75
$ws.Range("H136").Value = 21787942
$ws.Range("I136").Value = 45457850
$ws.Range("J136").Value = 90529.336
$ws.Range("K136").Value = 136373550
$ws.Range("L136").Value = 271588.008
$ws.Range("M136").Value = -136371000
$ws.Range("N136").Value = -276688.008

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2221.4211
$ws.Range("I94").Value = 1520.9
$ws.Range("K94").Value = 1520.9
$ws.Range("M94").Value = -1069.9
$ws.Range("H99").Value = 2289.25
$ws.Range("I99").Value = 1630.2222
$ws.Range("K99").Value = 1630.2222
$ws.Range("M99").Value = -132.2221999999999
$ws.Range("H107").Value = 2210
$ws.Range("I107").Value = 2224.0715
$ws.Range("K107").Value = 2224.0715
$ws.Range("M107").Value = -304.0715
$ws.Range("H134").Value = 32943.03
$ws.Range("I134").Value = 3354.8518
$ws.Range("J134").Value = 147068.86
$ws.Range("K134").Value = 10064.5554
$ws.Range("L134").Value = 441206.58
$ws.Range("M134").Value = -7529.555399999999
$ws.Range("N134").Value = -446276.58

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1072.45
$ws.Range("I16").Value = 778.625
$ws.Range("K16").Value = 778.625
$ws.Range("M16").Value = -491.625
$ws.Range("H113").Value = 1072.45
$ws.Range("I113").Value = 778.625
$ws.Range("K113").Value = 778.625
$ws.Range("M113").Value = 1391.375
$ws.Range("H122").Value = 5542.4165
$ws.Range("I122").Value = 3012
$ws.Range("K122").Value = 9036
$ws.Range("M122").Value = -6586

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 966
$ws.Range("I109").Value = 788.1667
$ws.Range("J109").Value = 3100
$ws.Range("K109").Value = 2364.5001
$ws.Range("L109").Value = 9300
$ws.Range("M109").Value = -1324.5001
$ws.Range("N109").Value = -11380

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1709.4667
$ws.Range("I97").Value = 1434
$ws.Range("K97").Value = 1434
$ws.Range("M97").Value = -938
$ws.Range("H122").Value = 38059.1
$ws.Range("I122").Value = 36732.332
$ws.Range("K122").Value = 110196.996
$ws.Range("M122").Value = -107746.996
$ws.Range("H126").Value = 4901
$ws.Range("I126").Value = 4830.2856
$ws.Range("K126").Value = 14490.8568
$ws.Range("M126").Value = -12020.8568
$ws.Range("H132").Value = 43481276
$ws.Range("I132").Value = 45456788
$ws.Range("J132").Value = 20000
$ws.Range("K132").Value = 136370364
$ws.Range("L132").Value = 60000
$ws.Range("M132").Value = -136367834
$ws.Range("N132").Value = -65060
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2830.9546
$ws.Range("I46").Value = 2623.4167
$ws.Range("J46").Value = 3080
$ws.Range("K46").Value = 2623.4167
$ws.Range("L46").Value = 3080
$ws.Range("M46").Value = -2435.4167
$ws.Range("N46").Value = -3456
$ws.Range("H93").Value = 250003440
$ws.Range("I93").Value = 333336260
$ws.Range("K93").Value = 333336260
$ws.Range("M93").Value = -333335012
$ws.Range("H100").Value = 3048.7058
$ws.Range("I100").Value = 3055.2
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 3055.2
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -2514.2
$ws.Range("N100").Value = -4082
$ws.Range("H111").Value = 58687.5
$ws.Range("I111").Value = 57988
$ws.Range("J111").Value = 59387
$ws.Range("K111").Value = 57988
$ws.Range("L111").Value = 59387
$ws.Range("M111").Value = -53898
$ws.Range("N111").Value = -67567
$ws.Range("H132").Value = 86290.125
$ws.Range("I132").Value = 50664.668
$ws.Range("K132").Value = 151994.004
$ws.Range("M132").Value = -149464.004

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3249
$ws.Range("I96").Value = 1662.6666
$ws.Range("J96").Value = 3777.7778
$ws.Range("K96").Value = 1662.6666
$ws.Range("L96").Value = 3777.7778
$ws.Range("M96").Value = -289.6666
$ws.Range("N96").Value = -6523.7778
$ws.Range("H122").Value = 7044.8335
$ws.Range("I122").Value = 4699.3335
$ws.Range("J122").Value = 9390.333000000001
$ws.Range("K122").Value = 14098.0005
$ws.Range("L122").Value = 28170.999
$ws.Range("M122").Value = -11648.0005
$ws.Range("N122").Value = -33070.999
$ws.Range("H126").Value = 5917.5405
$ws.Range("I126").Value = 5326.1724
$ws.Range("J126").Value = 8061.25
$ws.Range("K126").Value = 15978.5172
$ws.Range("L126").Value = 24183.75
$ws.Range("M126").Value = -13508.5172
$ws.Range("N126").Value = -29123.75
$ws.Range("H132").Value = 4086.125
$ws.Range("I132").Value = 4090.6155
$ws.Range("J132").Value = 4066.6667
$ws.Range("K132").Value = 12271.8465
$ws.Range("L132").Value = 12200.0001
$ws.Range("M132").Value = -9741.8465
$ws.Range("N132").Value = -17260.0001
